$wb = $excel.ActiveWorkbook

# Map of F-column updates (row -> new value) that apply identically to
# both the "展览" and "全部类型" sheets.
$updates = @{
    5  = 3063
    7  = 2290
    9  = 114
    11 = 1099
    13 = 48
    15 = 316
    16 = 274
    17 = 297
    21 = 53
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
